# Updated cryptos list - applies price/volume/coin changes per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-TextValue "D2" '64.300.15'
$ws.Range("E2").Value = '  -3.17%  '
Set-TextValue "D3" '3.170.97'
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue "D5" '564.28'
$ws.Range("E5").Value = '  -4.06%  '
Set-TextValue "D6" '170.37'
$ws.Range("E6").Value = '  -3.64%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.24%  '
Set-TextValue "D9" '3.168.62'
$ws.Range("E9").Value = '  -8.34%  '
$ws.Range("E10").Value = '  -6.79%  '
$ws.Range("E11").Value = '  -4.94%  '
Set-TextValue "D12" '0.396'
$ws.Range("E12").Value = '  -5.17%  '
Set-TextValue "D13" '3.717.64'
$ws.Range("E13").Value = '  -8.33%  '
$ws.Range("E14").Value = '  +0.87%  '
Set-TextValue "D15" '27.39'
$ws.Range("E15").Value = '  -7.46%  '
Set-TextValue "D16" '64.266.12'
$ws.Range("E16").Value = '  -3.02%  '
$ws.Range("E17").Value = '  -5.51%  '
Set-TextValue "D18" '3.169.90'
$ws.Range("E18").Value = '  -8.43%  '
$ws.Range("E19").Value = '  -3.95%  '
Set-TextValue "D20" '13.01'
$ws.Range("E20").Value = '  -5.75%  '
Set-TextValue "D21" '354.02'
$ws.Range("E21").Value = '  -5.45%  '
$ws.Range("E22").Value = '  -5.36%  '
$ws.Range("E23").Value = '  +0.03%  '
Set-TextValue "D24" '69.11'
$ws.Range("E24").Value = '  -5.83%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D25" '0.0000119'
$ws.Range("E25").Value = '  -5.81%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D26" '0.504'
$ws.Range("E26").Value = '  -6.22%  '
Set-TextValue "D27" '9.56'
$ws.Range("E27").Value = '  -3.53%  '
Set-TextValue "D28" '0.176'
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -4.05%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -5.40%  '
$ws.Range("E33").Value = '  -6.85%  '
$ws.Range("E34").Value = '  -5.85%  '
Set-TextValue "D36" '1.44'
$ws.Range("E36").Value = '  -7.57%  '
Set-TextValue "D37" '155.24'
$ws.Range("E37").Value = '  -3.84%  '
Set-TextValue "D38" '0.812'
$ws.Range("E38").Value = '  -8.14%  '
Set-TextValue "D39" '25.84'
$ws.Range("E39").Value = '  -9.41%  '
$ws.Range("E40").Value = '  -3.75%  '
$ws.Range("E41").Value = '  -6.45%  '
Set-TextValue "D42" '2.603.55'
$ws.Range("E42").Value = '  -6.04%  '
Set-TextValue "D43" '4.18'
$ws.Range("E43").Value = '  -7.21%  '
Set-TextValue "D44" '39.67'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D45" '0.0660'
$ws.Range("E45").Value = '  -4.74%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D46" '5.99'
$ws.Range("E46").Value = '  -7.15%  '
$ws.Range("E47").Value = '  -5.54%  '
Set-TextValue "D48" '321.48'
$ws.Range("E48").Value = '  -4.93%  '
$ws.Range("E49").Value = '  -7.29%  '
$ws.Range("E50").Value = '  -1.25%  '
$ws.Range("E51").Value = '  -0.04%  '
